# SF User list changes - 16 May - Initial
# Replace the sample user "Ashley Choi" with "Adarsh Patel" on the Users sheet,
# and leave the Users sheet as the active/selected sheet.

$wb = $excel.ActiveWorkbook

$usersSheet = $wb.Worksheets.Item("Users")
$usersSheet.Range("A2").Value = "Adarsh Patel"

# Make the Users sheet the active sheet/tab.
$usersSheet.Activate() | Out-Null
$usersSheet.Select() | Out-Null
$usersSheet.Range("A2").Select() | Out-Null
